$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.23145318031311
$ws.Range("B1").Value = 2.502643346786499
$ws.Range("C1").Value = 4.336839199066162
$ws.Range("D1").Value = 2.54498291015625
$ws.Range("E1").Value = 1.078810214996338
